$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "59.699.93"
$ws.Cells.Item(2, 5).Value = "  +0.33%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.648.98"
$ws.Cells.Item(3, 5).Value = "  +1.45%  "
$ws.Cells.Item(4, 5).Value = "  +0.01%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "538.35"
$ws.Cells.Item(5, 5).Value = "  +0.28%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "145.69"
$ws.Cells.Item(6, 5).Value = "  +3.00%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.999"
$ws.Cells.Item(7, 5).Value = "  +0.04%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.574"
$ws.Cells.Item(8, 5).Value = "  +0.82%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "2.666.08"
$ws.Cells.Item(9, 5).Value = "  +1.59%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "6.65"
$ws.Cells.Item(10, 5).Value = "  +2.36%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.103"
$ws.Cells.Item(11, 5).Value = "  -0.15%  "
$ws.Cells.Item(12, 5).Value = "  +0.29%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.134"
$ws.Cells.Item(13, 5).Value = "  -0.81%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "3.122.65"
$ws.Cells.Item(14, 5).Value = "  +1.75%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "59.669.60"
$ws.Cells.Item(15, 5).Value = "  +0.41%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "21.25"
$ws.Cells.Item(16, 5).Value = "  +3.01%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.663.68"
$ws.Cells.Item(17, 5).Value = "  +1.43%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.0000135"
$ws.Cells.Item(18, 5).Value = "  +0.95%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "342.10"
$ws.Cells.Item(19, 5).Value = "  -1.22%  "
$ws.Cells.Item(20, 5).Value = "  +1.36%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "10.43"
$ws.Cells.Item(21, 5).Value = "  +2.69%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.32"
$ws.Cells.Item(22, 5).Value = "  -1.05%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.998"
$ws.Cells.Item(23, 5).Value = "  -0.04%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "66.77"
$ws.Cells.Item(24, 5).Value = "  -0.54%  "
$ws.Cells.Item(25, 5).Value = "  +1.96%  "
$ws.Cells.Item(26, 5).Value = "  -1.08%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.999"
$ws.Cells.Item(27, 5).Value = "  +0.01%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.33"
$ws.Cells.Item(28, 5).Value = "  +1.39%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.0₃0754"
$ws.Cells.Item(29, 5).Value = "  +0.25%  "
$ws.Cells.Item(30, 5).Value = "  -0.02%  "
$ws.Cells.Item(31, 5).Value = "  +1.06%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "5.86"
$ws.Cells.Item(32, 5).Value = "  -0.42%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "18.99"
$ws.Cells.Item(33, 5).Value = "  +0.22%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "150.98"
$ws.Cells.Item(34, 5).Value = "  +1.17%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.03"
$ws.Cells.Item(35, 5).Value = "  +0.30%  "
$ws.Cells.Item(36, 5).Value = "  +1.06%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.843"
$ws.Cells.Item(37, 5).Value = "  -0.51%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.841"
$ws.Cells.Item(38, 5).Value = "  -0.17%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.46"
$ws.Cells.Item(39, 5).Value = "  -0.71%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "291.59"
$ws.Cells.Item(40, 5).Value = "  +4.77%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "3.60"
$ws.Cells.Item(41, 5).Value = "  +1.30%  "
$ws.Cells.Item(42, 5).Value = "  +0.09%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.607"
$ws.Cells.Item(43, 5).Value = "  +0.82%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "19.44"
$ws.Cells.Item(44, 5).Value = "  +3.51%  "
$ws.Cells.Item(45, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "10.72"
$ws.Cells.Item(45, 5).Value = "  -0.38%  "
$ws.Cells.Item(46, 2).Value = "Hedera"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0538"
$ws.Cells.Item(46, 5).Value = "  +2.45%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0949"
$ws.Cells.Item(47, 5).Value = "  -1.40%  "
$ws.Cells.Item(48, 2).Value = "VeChain"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0228"
$ws.Cells.Item(48, 5).Value = "  +1.68%  "
$ws.Cells.Item(49, 2).Value = "Maker"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.980.16"
$ws.Cells.Item(49, 5).Value = "  +1.27%  "
$ws.Cells.Item(50, 2).Value = "RenderToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "4.63"
$ws.Cells.Item(50, 5).Value = "  +1.97%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "18.41"
$ws.Cells.Item(51, 5).Value = "  -0.10%  "
